$d = $word.ActiveDocument

# --- Paragraph 1: insert a missing space between "proposed" and "NRT" ---
$d.Content.Find.Execute("proposedNRT program.", $true, $false, $false, $false, $false, $true, 1, $false, "proposed NRT program.", 2)

# --- Paragraph 3: replace the "1 page" placeholder with the full training narrative ---
$d.Content.Find.Execute("1 page", $false, $false, $false, $false, $false, $true, 1, $false, "Over the past five years, PI O’Meara has had four graduate students in his lab. All are still currently enrolled: two are on schedule to receive their PhD in ecology and evolutionary biology as well as a Masters in statistics this semester, another recently received a DDIG award and is on track for graduation on schedule, and a fourth recently took his PhD qualifying exam. The two students planning to finish this semester have decided to pursue careers outside of academia, and enrolled in and successfully completed a program to earn a Masters in statistics while in a PhD program with this intention. They both had internships at the Tennessee Valley Authority (one received an offer of a job once she graduated), and one has also interned with our athletic department analyzing academic progress of athletes. Half of O’Meara’s students identify as women, and one identifies as Hispanic. O’Meara also serves on approximately one-third of graduate student committees in the Ecology and Evolutionary Biology department and has also served on student committees in Entomology, Earth & Planetary Sciences, Microbiology, and Genome Sciences and Technology. O’Meara has also served on EEB’s graduate admission committee and now is associate head for graduate affairs in the department; as part of this, he has run training for graduate students in grant writing. As associate director for postdoctoral training for NIMBioS, he has also organized training sessions for postdocs pursuing careers in biology, math, and statistics.", 2)

# --- Restore the _GoBack bookmark at its original (mid-sentence) position,
#     right after "...As associate director for postdoctoral training " and
#     before "for NIMBioS, he has also organized..." ---
$p3 = $d.Paragraphs(3).Range
$bmPos = $p3.Start + 1411
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
